# Applies updated sensitivity/calculus results to the daily_model workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Model Accuracy (-0.25, 0.25, 0.25)"
# Add new columns (Market threshold, Market min, Market max, Recall, Precision)
# and refresh the Accuracy (%) values.
# ---------------------------------------------------------------------------
$wsAccuracy = $wb.Worksheets.Item("Model Accuracy (-0.25, 0.25, 0.25)")

# Replicate the bold/bordered header formatting used by B1 onto the new
# header cells before filling in their text.
$wsAccuracy.Range("B1").Copy()
$wsAccuracy.Range("C1:G1").PasteSpecial(-4122)

$wsAccuracy.Range("C1").Value = "Market threshold"
$wsAccuracy.Range("D1").Value = "Market min"
$wsAccuracy.Range("E1").Value = "Market max"
$wsAccuracy.Range("F1").Value = "Recall"
$wsAccuracy.Range("G1").Value = "Precision"

# TOTALENERGIES SE
$wsAccuracy.Range("B2").Value = 54.76772616136919
$wsAccuracy.Range("C2").Value = 0.05450546436368681
$wsAccuracy.Range("D2").Value = -15.55441
$wsAccuracy.Range("E2").Value = 15.06418
$wsAccuracy.Range("F2").Value = 11.11111111111111
$wsAccuracy.Range("G2").Value = 1.923076923076923

# FMC CORP
$wsAccuracy.Range("B3").Value = 32.09046454767726
$wsAccuracy.Range("C3").Value = 0.009583939973006913
$wsAccuracy.Range("D3").Value = -19.35264
$wsAccuracy.Range("E3").Value = 13.70093
$wsAccuracy.Range("F3").Value = 8.042895442359249
$wsAccuracy.Range("G3").Value = 22.90076335877863

# BP PLC
$wsAccuracy.Range("B4").Value = 81.17359413202934
$wsAccuracy.Range("C4").Value = 0.04158117063764853
$wsAccuracy.Range("D4").Value = -18.75314
$wsAccuracy.Range("E4").Value = 23.33066
$wsAccuracy.Range("F4").Value = 0
$wsAccuracy.Range("G4").Value = 0

# STORA ENSO
$wsAccuracy.Range("B5").Value = 69.00977995110024
$wsAccuracy.Range("C5").Value = 0.02983403801513819
$wsAccuracy.Range("D5").Value = -12.78028
$wsAccuracy.Range("E5").Value = 12.42348
$wsAccuracy.Range("F5").Value = 3.636363636363636
$wsAccuracy.Range("G5").Value = 5.88235294117647

# BHP GROUP
$wsAccuracy.Range("B6").Value = 81.72371638141809
$wsAccuracy.Range("C6").Value = 0.08368817696170747
$wsAccuracy.Range("D6").Value = -16.47904
$wsAccuracy.Range("E6").Value = 14.94325
$wsAccuracy.Range("F6").Value = 0
$wsAccuracy.Range("G6").Value = 0

# ---------------------------------------------------------------------------
# Sheet: "Confusion Matrix TOTALENERGIES SE (-0.25, 0.25, 0.25)"
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("Confusion Matrix TOTALENERGIES SE (-0.25, 0.25, 0.25)")
$wsTotal.Range("B3").Value = 6
$wsTotal.Range("C3").Value = 889
$wsTotal.Range("D3").Value = 4

# ---------------------------------------------------------------------------
# Sheet: "Confusion Matrix FMC CORP (-0.25, 0.25, 0.25)"
# ---------------------------------------------------------------------------
$wsFmc = $wb.Worksheets.Item("Confusion Matrix FMC CORP (-0.25, 0.25, 0.25)")
$wsFmc.Range("B2").Value = 30
$wsFmc.Range("C2").Value = 70
$wsFmc.Range("D2").Value = 31
$wsFmc.Range("B3").Value = 210
$wsFmc.Range("C3").Value = 367
$wsFmc.Range("D3").Value = 197
$wsFmc.Range("B4").Value = 133
$wsFmc.Range("C4").Value = 215
$wsFmc.Range("D4").Value = 128

# ---------------------------------------------------------------------------
# Sheet: "Confusion Matrix BP PLC (-0.25, 0.25, 0.25)"
# ---------------------------------------------------------------------------
$wsBp = $wb.Worksheets.Item("Confusion Matrix BP PLC (-0.25, 0.25, 0.25)")
$wsBp.Range("B3").Value = 33
$wsBp.Range("C3").Value = 1318
$wsBp.Range("D3").Value = 32
$wsBp.Range("B4").Value = 7
$wsBp.Range("C4").Value = 196
$wsBp.Range("D4").Value = 10

# ---------------------------------------------------------------------------
# Sheet: "Confusion Matrix STORA ENSO (-0.25, 0.25, 0.25)"
# ---------------------------------------------------------------------------
$wsStora = $wb.Worksheets.Item("Confusion Matrix STORA ENSO (-0.25, 0.25, 0.25)")
$wsStora.Range("B2").Value = 4
$wsStora.Range("C2").Value = 63
$wsStora.Range("B3").Value = 88
$wsStora.Range("C3").Value = 1115
$wsStora.Range("D3").Value = 96
$wsStora.Range("B4").Value = 18
$wsStora.Range("C4").Value = 180
$wsStora.Range("D4").Value = 10

# ---------------------------------------------------------------------------
# Sheet: "Confusion Matrix BHP GROUP (-0.25, 0.25, 0.25)"
# ---------------------------------------------------------------------------
$wsBhp = $wb.Worksheets.Item("Confusion Matrix BHP GROUP (-0.25, 0.25, 0.25)")
$wsBhp.Range("B2").Value = 0
$wsBhp.Range("C2").Value = 158
$wsBhp.Range("D2").Value = 1
$wsBhp.Range("B3").Value = 4
$wsBhp.Range("C3").Value = 1337
$wsBhp.Range("D3").Value = 2
$wsBhp.Range("B4").Value = 0
$wsBhp.Range("C4").Value = 78
